$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 16966
$ws.Range("E2").Value = 857
$ws.Range("F2").Value = 879
$ws.Range("G2").Value = 240
$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 95
$ws.Range("J2").Value = -1
$ws.Range("K2").Value = 11593
$ws.Range("L2").Value = 9685
$ws.Range("M2").Value = 1907
$ws.Range("N2").Value = 1922
$ws.Range("O2").Value = -14
$ws.Range("P2").Value = 323
$ws.Range("Q2").Value = 1025
$ws.Range("R2").Value = -1097
$ws.Range("S2").Value = 214
$ws.Range("T2").Value = 755
$ws.Range("U2").Value = 270
$ws.Range("V2").Value = 6944
$ws.Range("W2").Value = 5.05
$ws.Range("X2").Value = 0.5600000000000001
$ws.Range("Y2").Value = 4.92
$ws.Range("Z2").Value = 0.83
$ws.Range("AA2").Value = 507.76
$ws.Range("AB2").Value = 506.4
$ws.Range("AC2").Value = 148
$ws.Range("AD2").Value = 20.53
$ws.Range("AE2").Value = 3227
$ws.Range("AF2").Value = 0.9399999999999999
$ws.Range("AG2").Value = 25
$ws.Range("AH2").Value = 0.83
$ws.Range("AI2").Value = 15.63
$ws.Range("AJ2").Value = 64561210

# Row 3
$ws.Range("D3").Value = 16795
$ws.Range("E3").Value = 1100
$ws.Range("F3").Value = 1100
$ws.Range("G3").Value = 832
$ws.Range("H3").Value = 770
$ws.Range("I3").Value = 769
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 12995
$ws.Range("L3").Value = 10307
$ws.Range("M3").Value = 2688
$ws.Range("N3").Value = 2686
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 323
$ws.Range("Q3").Value = 1453
$ws.Range("R3").Value = -1827
$ws.Range("S3").Value = 339
$ws.Range("T3").Value = 732
$ws.Range("U3").Value = 722
$ws.Range("V3").Value = 7282
$ws.Range("W3").Value = 6.55
$ws.Range("X3").Value = 4.58
$ws.Range("Y3").Value = 33.38
$ws.Range("Z3").Value = 6.26
$ws.Range("AA3").Value = 383.51
$ws.Range("AB3").Value = 724.6
$ws.Range("AC3").Value = 1191
$ws.Range("AD3").Value = 3.26
$ws.Range("AE3").Value = 4510
$ws.Range("AF3").Value = 0.86
$ws.Range("AG3").Value = 50
$ws.Range("AH3").Value = 1.29
$ws.Range("AI3").Value = 3.87
$ws.Range("AJ3").Value = 64561210

# Row 4
$ws.Range("D4").Value = 16161
$ws.Range("E4").Value = 1053
$ws.Range("F4").Value = 1053
$ws.Range("G4").Value = 843
$ws.Range("H4").Value = 588
$ws.Range("I4").Value = 583
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 13398
$ws.Range("L4").Value = 10146
$ws.Range("M4").Value = 3253
$ws.Range("N4").Value = 3246
$ws.Range("O4").Value = 7
$ws.Range("P4").Value = 323
$ws.Range("Q4").Value = 773
$ws.Range("R4").Value = -50
$ws.Range("S4").Value = -360
$ws.Range("T4").Value = 532
$ws.Range("U4").Value = 241
$ws.Range("V4").Value = 7046
$ws.Range("W4").Value = 6.52
$ws.Range("X4").Value = 3.64
$ws.Range("Y4").Value = 19.64
$ws.Range("Z4").Value = 4.46
$ws.Range("AA4").Value = 311.93
$ws.Range("AB4").Value = 888.04
$ws.Range("AC4").Value = 902
$ws.Range("AD4").Value = 3.72
$ws.Range("AE4").Value = 5450
$ws.Range("AF4").Value = 0.62
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 1.49
$ws.Range("AI4").Value = 5.11
$ws.Range("AJ4").Value = 64561210

# Row 5
$ws.Range("D5").Value = 14434
$ws.Range("E5").Value = 358
$ws.Range("F5").Value = 358
$ws.Range("G5").Value = -408
$ws.Range("H5").Value = -418
$ws.Range("I5").Value = -423
$ws.Range("J5").Value = 5
$ws.Range("K5").Value = 12742
$ws.Range("L5").Value = 10039
$ws.Range("M5").Value = 2703
$ws.Range("N5").Value = 2694
$ws.Range("O5").Value = 9
$ws.Range("P5").Value = 323
$ws.Range("Q5").Value = 594
$ws.Range("R5").Value = -392
$ws.Range("S5").Value = -44
$ws.Range("T5").Value = 856
$ws.Range("U5").Value = -262
$ws.Range("V5").Value = 6959
$ws.Range("W5").Value = 2.48
$ws.Range("X5").Value = -2.9
$ws.Range("Y5").Value = -14.23
$ws.Range("Z5").Value = -3.2
$ws.Range("AA5").Value = 371.44
$ws.Range("AB5").Value = 759.73
$ws.Range("AC5").Value = -655
$ws.Range("AD5").Value = -4
$ws.Range("AE5").Value = 4397
$ws.Range("AF5").Value = 0.6
$ws.Range("AG5").Value = 25
$ws.Range("AH5").Value = 0.95
$ws.Range("AI5").Value = -3.62
$ws.Range("AJ5").Value = 64561210

# Row 6
$ws.Range("D6").Value = 15281
$ws.Range("E6").Value = 539
$ws.Range("F6").Value = 539
$ws.Range("G6").Value = -305
$ws.Range("H6").Value = -332
$ws.Range("I6").Value = -337
$ws.Range("K6").Value = 12118
$ws.Range("L6").Value = 9933
$ws.Range("M6").Value = 2185
$ws.Range("N6").Value = 2176
$ws.Range("P6").Value = 323
$ws.Range("Q6").Value = 192
$ws.Range("R6").Value = -215
$ws.Range("S6").Value = 40
$ws.Range("T6").Value = 524
$ws.Range("U6").Value = -332
$ws.Range("V6").Value = 7018
$ws.Range("W6").Value = 3.53
$ws.Range("X6").Value = -2.17
$ws.Range("Y6").Value = -13.82
$ws.Range("Z6").Value = -2.67
$ws.Range("AA6").Value = 454.62
$ws.Range("AB6").Value = 654.5
$ws.Range("AC6").Value = -521
$ws.Range("AD6").Value = -3.95
$ws.Range("AE6").Value = 3551
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6").Value = 25
$ws.Range("AH6").Value = 1.21
$ws.Range("AI6").Value = -4.55
$ws.Range("AJ6").Value = 64561210

# Clear cells in rows 7-9 (removed in diff)
# Row 7
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

